# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 77
$ws1.Range("F6").Value = 538
$ws1.Range("F7").Value = 1652
$ws1.Range("F8").Value = 8
$ws1.Range("F11").Value = 1537
$ws1.Range("F13").Value = 49
$ws1.Range("F15").Value = 254
$ws1.Range("F21").Value = 274
$ws1.Range("F24").Value = 213

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 77
$ws4.Range("F6").Value = 538
$ws4.Range("F7").Value = 1652
$ws4.Range("F9").Value = 8
$ws4.Range("F12").Value = 1537
$ws4.Range("F14").Value = 49
$ws4.Range("F16").Value = 254
$ws4.Range("F22").Value = 274
$ws4.Range("F25").Value = 213
